# Corrected excel sheets for application fix issues
# Applies the numeric corrections (and resulting selection moves) that the
# author made to the repayment-schedule / summary / transactions sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()

$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 211.3
$wsSummary.Range("E3").Value = 114.31

$wsSummary.Range("D4").Select()

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Activate()

# Blank helper cell in row 2 slides from column P to column O.
$wsSchedule.Range("P2").Copy()
$wsSchedule.Range("O2").PasteSpecial(-4122)  # xlPasteFormats
$wsSchedule.Range("P2").Clear()
$wsSchedule.Application.CutCopyMode = $false

# Row 4 (installment #2)
$wsSchedule.Range("B4").Value = 31
$wsSchedule.Range("C4").Value = 42095
$wsSchedule.Range("H4").Value = 41.69
$wsSchedule.Range("K4").Value = 950.78
$wsSchedule.Range("P4").Value = 950.78

# Row 5 (installment #3)
$wsSchedule.Range("B5").Value = 30
$wsSchedule.Range("C5").Value = 42125
$wsSchedule.Range("H5").Value = 31.38
$wsSchedule.Range("K5").Value = 940.47
$wsSchedule.Range("P5").Value = 940.47

# Row 6 (installment #4)
$wsSchedule.Range("B6").Value = 31
$wsSchedule.Range("C6").Value = 42156
$wsSchedule.Range("H6").Value = 23.16
$wsSchedule.Range("K6").Value = 932.25
$wsSchedule.Range("P6").Value = 932.25

# Row 7 (installment #5)
$wsSchedule.Range("B7").Value = 30
$wsSchedule.Range("C7").Value = 42186
$wsSchedule.Range("H7").Value = 13.45
$wsSchedule.Range("K7").Value = 922.54
$wsSchedule.Range("P7").Value = 922.54

# Row 8 (installment #6)
$wsSchedule.Range("B8").Value = 31
$wsSchedule.Range("C8").Value = 42217
$wsSchedule.Range("H8").Value = 4.63
$wsSchedule.Range("K8").Value = 459.18
$wsSchedule.Range("P8").Value = 459.18

# Column L widened slightly (best-fit recalculated after the edits above).
$wsSchedule.Columns.Item(12).ColumnWidth = 7.3

$wsSchedule.Range("M7").Select()

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()

$wsTrans.Range("A2").Value = 6496
$wsTrans.Range("A3").Value = 6494

$wsTrans.Range("D3").Select()
